$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Fix the title: "OBSERVACIONES DEL LA PRACTICA"
#    -> "OBSERVACIONES DE LA PRÁCTICA"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "OBSERVACIONES DEL LA PRACTICA", $true, $false, $false, $false, $false,
    $true, 1, $false, "OBSERVACIONES DE LA PRÁCTICA", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Add a third "Estudiante" line (right after "Estudiante 2 Cod XXXX"),
#    matching the formatting of the existing student lines.
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("Estudiante 2 Cod XXXX", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $stuParagraph = $d.Paragraphs.Item(3)
    $stuParagraph.Range.InsertParagraphAfter()
    $newPar = $d.Paragraphs.Item(4)
    $newPar.Range.Text = "Estudiante 3 Cod XXXX"
}

# ---------------------------------------------------------------------
# 3) Move the inverted question mark in question 2 so the sentence reads
#    "...en un BST, ¿cree que..." instead of starting with "¿Si tuviera".
# ---------------------------------------------------------------------
$oldQ2 = "¿Si tuviera que responder esa misma consulta y la información estuviera en tablas de hash y no en un BST, cree que el tiempo de respuesta sería mayor o menor? ¿Por qué?"
$newQ2 = "Si tuviera que responder esa misma consulta y la información estuviera en tablas de hash y no en un BST, ¿cree que el tiempo de respuesta sería mayor o menor? ¿Por qué?"

$d.Content.Find.Execute($oldQ2, $true, $false, $false, $false, $false, $true, 1, $false, $newQ2, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) The empty paragraph right under the "Preguntas de análisis" heading
#    carried a redundant explicit "Normal" paragraph style; clear it so
#    the style is implicit (default) again, like the rest of the body.
# ---------------------------------------------------------------------
$emptyPar = $d.Paragraphs.Item(6)
if ($emptyPar.Range.Text -eq [char]13 -and $emptyPar.Style.NameLocal -eq "Normal") {
    $emptyPar.Range.Style = $d.Styles.Item("Normal")
    $emptyPar.Range.LanguageID = "es-CO"
    $emptyPar.Range.NoProofing = 0
}
